$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("Test Cases")

# Set Runmode column (C) to "Y" for all test case rows, and Results (D) to "PASS"
$wsCases.Range("C2").Value = "Y"
$wsCases.Range("D2").Value = "PASS"

$wsCases.Range("C3").Value = "Y"
$wsCases.Range("D3").Value = "PASS"

$wsCases.Range("C4").Value = "Y"
$wsCases.Range("D4").Value = "PASS"

# Scroll the "Test Cases" sheet so column C is the top-left visible cell
$wsCases.Application.ActiveWindow.ScrollColumn = 3

$wb.Save()
